$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("MTF_UT_0020", "test_add_child"),
    @("MTF_UT_0021", "test_add_spouse"),
    @("MTF_UT_0022", "test_get_relationship"),
    @("MFT_IT_0003", "test_add_child"),
    @("MFT_IT_0004", "test_add_spouse"),
    @("MFT_IT_0005", "test_get_relationship")
)

$row = 25
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

$ws.Range("I25").Select()
